$d = $word.ActiveDocument

# Update the date/day title
$d.Content.Find.Execute("2024-09-11 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-09-12 Thursday", 2)

# Update the division problems in the table (positional, since several
# values repeat and are not unique across the document)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "23÷8="
$t.Cell(1,2).Range.Text  = "58÷5="
$t.Cell(1,3).Range.Text  = "43÷8="
$t.Cell(1,4).Range.Text  = "58÷2="
$t.Cell(1,5).Range.Text  = "70÷5="

$t.Cell(5,1).Range.Text  = "19÷8="
$t.Cell(5,2).Range.Text  = "36÷9="
$t.Cell(5,3).Range.Text  = "73÷6="
$t.Cell(5,4).Range.Text  = "75÷5="
$t.Cell(5,5).Range.Text  = "76÷3="

$t.Cell(9,1).Range.Text  = "88÷7="
$t.Cell(9,2).Range.Text  = "66÷8="
$t.Cell(9,3).Range.Text  = "82÷6="
$t.Cell(9,4).Range.Text  = "58÷4="
$t.Cell(9,5).Range.Text  = "83÷3="

$t.Cell(13,1).Range.Text = "24÷6="
$t.Cell(13,2).Range.Text = "72÷9="
$t.Cell(13,3).Range.Text = "37÷4="
$t.Cell(13,4).Range.Text = "81÷9="
$t.Cell(13,5).Range.Text = "36÷8="

$t.Cell(17,2).Range.Text = "58÷4="
$t.Cell(17,3).Range.Text = "60÷5="
$t.Cell(17,4).Range.Text = "73÷9="
$t.Cell(17,5).Range.Text = "38÷8="
